$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore the workbook color palette to Excel defaults (as happens when the
# workbook is next saved by Excel after being produced by a non-Excel tool).
$wb.ResetColors()

# Populate the newly introduced "Circumference" column (C) with the results
# of the updated circumference algorithm, for every data row (2-94).
$ws.Cells.Item(2, 3).Value = 1119.134118676186
$ws.Cells.Item(3, 3).Value = 914.9747383594513
$ws.Cells.Item(4, 3).Value = 506.6589421033859
$ws.Cells.Item(5, 3).Value = 288.2497808933258
$ws.Cells.Item(6, 3).Value = 724.7493426799774
$ws.Cells.Item(7, 3).Value = 406.7594473361969
$ws.Cells.Item(8, 3).Value = 845.0752435922623
$ws.Cells.Item(9, 3).Value = 701.3523740768433
$ws.Cells.Item(10, 3).Value = 1160.721993088722
$ws.Cells.Item(11, 3).Value = 1213.567665100098
$ws.Cells.Item(12, 3).Value = 1619.498685359955
$ws.Cells.Item(13, 3).Value = 703.6122596263885
$ws.Cells.Item(14, 3).Value = 286.1492756605148
$ws.Cells.Item(15, 3).Value = 173.3797236680984
$ws.Cells.Item(16, 3).Value = 195.5807341337204
$ws.Cells.Item(17, 3).Value = 1430.805253386497
$ws.Cells.Item(18, 3).Value = 452.558436870575
$ws.Cells.Item(19, 3).Value = 216.3502861261368
$ws.Cells.Item(20, 3).Value = 704.4234417676926
$ws.Cells.Item(21, 3).Value = 242.2497808933258
$ws.Cells.Item(22, 3).Value = 650.7665876150131
$ws.Cells.Item(23, 3).Value = 1221.248904466629
$ws.Cells.Item(24, 3).Value = 306.433546423912
$ws.Cells.Item(25, 3).Value = 341.2792184352875
$ws.Cells.Item(26, 3).Value = 282.3919162750244
$ws.Cells.Item(27, 3).Value = 291.9482651948929
$ws.Cells.Item(28, 3).Value = 164.4091612100601
$ws.Cells.Item(29, 3).Value = 392.0315254926682
$ws.Cells.Item(30, 3).Value = 295.3208485841751
$ws.Cells.Item(31, 3).Value = 206.1076455116272
$ws.Cells.Item(32, 3).Value = 116.5685415267944
$ws.Cells.Item(33, 3).Value = 2429.586121559143
$ws.Cells.Item(34, 3).Value = 100.5685415267944
$ws.Cells.Item(35, 3).Value = 435.2447285652161
$ws.Cells.Item(36, 3).Value = 439.7300097942352
$ws.Cells.Item(37, 3).Value = 110.811182141304
$ws.Cells.Item(38, 3).Value = 267.3208485841751
$ws.Cells.Item(39, 3).Value = 60.76955199241638
$ws.Cells.Item(40, 3).Value = 362.7178171873093
$ws.Cells.Item(41, 3).Value = 936.3717069625854
$ws.Cells.Item(42, 3).Value = 395.7888848781586
$ws.Cells.Item(43, 3).Value = 305.3624787330627
$ws.Cells.Item(44, 3).Value = 365.64674949646
$ws.Cells.Item(45, 3).Value = 36.38477599620819
$ws.Cells.Item(46, 3).Value = 254.534051656723
$ws.Cells.Item(47, 3).Value = 418.1736608743668
$ws.Cells.Item(48, 3).Value = 295.3624787330627
$ws.Cells.Item(49, 3).Value = 1075.692555546761
$ws.Cells.Item(50, 3).Value = 496.6000670194626
$ws.Cells.Item(51, 3).Value = 648.3990565538406
$ws.Cells.Item(52, 3).Value = 238.4924215078354
$ws.Cells.Item(53, 3).Value = 675.5950146913528
$ws.Cells.Item(54, 3).Value = 407.7300097942352
$ws.Cells.Item(55, 3).Value = 627.068103313446
$ws.Cells.Item(56, 3).Value = 376.6589421033859
$ws.Cells.Item(57, 3).Value = 257.2203433513641
$ws.Cells.Item(58, 3).Value = 468.558436870575
$ws.Cells.Item(59, 3).Value = 392.2741661071777
$ws.Cells.Item(60, 3).Value = 243.5634891986847
$ws.Cells.Item(61, 3).Value = 230.4507913589478
$ws.Cells.Item(62, 3).Value = 722.4234417676926
$ws.Cells.Item(63, 3).Value = 216.2081507444382
$ws.Cells.Item(64, 3).Value = 1147.325024485588
$ws.Cells.Item(65, 3).Value = 1266.219466924667
$ws.Cells.Item(66, 3).Value = 397.7888848781586
$ws.Cells.Item(67, 3).Value = 187.9655101299286
$ws.Cells.Item(68, 3).Value = 1392.101716756821
$ws.Cells.Item(69, 3).Value = 302.4924215078354
$ws.Cells.Item(70, 3).Value = 819.6193999052048
$ws.Cells.Item(71, 3).Value = 289.9066350460052
$ws.Cells.Item(72, 3).Value = 319.2619735002518
$ws.Cells.Item(73, 3).Value = 249.1787132024765
$ws.Cells.Item(74, 3).Value = 551.9554054737091
$ws.Cells.Item(75, 3).Value = 234.1076455116272
$ws.Cells.Item(76, 3).Value = 357.4041088819504
$ws.Cells.Item(77, 3).Value = 994.8986183404922
$ws.Cells.Item(78, 3).Value = 344.7766922712326
$ws.Cells.Item(79, 3).Value = 499.6711347103119
$ws.Cells.Item(80, 3).Value = 349.2030984163284
$ws.Cells.Item(81, 3).Value = 327.2619735002518
$ws.Cells.Item(82, 3).Value = 466.1736608743668
$ws.Cells.Item(83, 3).Value = 387.8650048971176
$ws.Cells.Item(84, 3).Value = 624.4406867027283
$ws.Cells.Item(85, 3).Value = 822.4478269815445
$ws.Cells.Item(86, 3).Value = 44.62741661071777
$ws.Cells.Item(87, 3).Value = 238.6934319734573
$ws.Cells.Item(88, 3).Value = 372.1736608743668
$ws.Cells.Item(89, 3).Value = 312.0487704277039
$ws.Cells.Item(90, 3).Value = 235.3797236680984
$ws.Cells.Item(91, 3).Value = 280.8771975040436
$ws.Cells.Item(92, 3).Value = 257.4213538169861
$ws.Cells.Item(93, 3).Value = 207.82337474823
$ws.Cells.Item(94, 3).Value = 123.6396092176437

